# "before elimination of J7" -- rescale the R/C component values in row 2
# (each literal is shifted one decade; the dependent formula cells in
# C/E/G/I/K/L/M recalc automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 59000
$ws.Range("B2").Value = 0.0000000033
$ws.Range("D2").Value = 0.000000000018
$ws.Range("F2").Value = 0.0000000082
$ws.Range("H2").Value = 348
$ws.Range("J2").Value = 17700

# Scroll the sheet view so column G is left-most visible (was column C).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
